$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.280.25"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").Value = "'3.384.35"
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'587.30"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("D6").Value = "'179.18"
$ws.Range("E6").Value = "  +1.32%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'0.596"
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("E9").Value = "  +5.53%  "
$ws.Range("D10").Value = "'0.590"
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").Value = "'48.37"
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("D12").Value = "'0.0000281"
$ws.Range("E12").Value = "  +2.97%  "
$ws.Range("D13").Value = "'680.23"
$ws.Range("E13").Value = "  -2.79%  "
$ws.Range("D14").Value = "'8.60"
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("D15").Value = "'3.929.16"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").Value = "'69.366.64"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "'3.415.37"
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("D19").Value = "'17.64"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("D21").Value = "'0.904"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").Value = "'5.41"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'17.19"
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("D24").Value = "'103.36"
$ws.Range("E24").Value = "  +3.37%  "
$ws.Range("D25").Value = "'3.93"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").Value = "'9.63"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").Value = "'33.92"
$ws.Range("E28").Value = "  +2.76%  "
$ws.Range("D29").Value = "'8.71"
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").Value = "'6.95"
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("D31").Value = "'562.86"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "'11.12"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").Value = "'3.53"
$ws.Range("E34").Value = "  +4.45%  "
$ws.Range("D35").Value = "'58.59"
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "'3.685.84"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'35.71"
$ws.Range("E38").Value = "  +2.85%  "
$ws.Range("E39").Value = "  +4.30%  "
$ws.Range("D40").Value = "'3.25"
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("D41").Value = "'2.68"
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("D42").Value = "'0.0₃0696"
$ws.Range("E42").Value = "  +3.24%  "
$ws.Range("D43").Value = "'0.339"
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("D44").Value = "'0.0423"
$ws.Range("E44").Value = "  +3.65%  "
$ws.Range("D45").Value = "'3.29"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").Value = "'2.67"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("D48").Value = "'1.42"
$ws.Range("E48").Value = "  +5.31%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'133.24"
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("D51").Value = "'2.61"
$ws.Range("E51").Value = "  +3.20%  "
